$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (reflects the "through" date moving from 12-19 to 12-20)
$ws.Name = "Through 2022-12-20"

# Update the row label for December
$ws.Range("A13").Value = "December (through 12-20)"

# Updated December figures (row 13)
$ws.Range("B13").Value = 27
$ws.Range("C13").Value = 64
$ws.Range("D13").Value = 77
$ws.Range("E13").Value = 44
$ws.Range("F13").Value = 33
$ws.Range("G13").Value = 93
$ws.Range("H13").Value = 140
$ws.Range("I13").Value = 84

# Updated Total figures (row 14)
$ws.Range("B14").Value = 318
$ws.Range("C14").Value = 627
$ws.Range("D14").Value = 898
$ws.Range("E14").Value = 726
$ws.Range("F14").Value = 567
$ws.Range("G14").Value = 1357
$ws.Range("H14").Value = 1783
$ws.Range("I14").Value = 1601
